# "break out stock.yaml completed"
#
# 1) On the "day" sheet, rows 85-89 (Coromandel, GNFC, Vedanta, Nmdc,
#    Aditya Birla Capital) had their bsecode (column D) stored as text;
#    normalize them to numbers like every other row in the column.
# 2) On the "week" sheet, append three newly scraped rows (94-96:
#    Colgate Palmolive, United Breweries, National Aluminium) that were
#    broken out from stock.yaml, extending the sheet from A1:I93 to
#    A1:I96. Their bsecode values keep the as-scraped text formatting
#    (matching how the day-sheet rows originally looked before cleanup).

$wb = $excel.ActiveWorkbook

$dayWs = $wb.Worksheets.Item("day")
$dayWs.Range("D85").Value = 506395
$dayWs.Range("D86").Value = 500670
$dayWs.Range("D87").Value = 500295
$dayWs.Range("D88").Value = 526371
$dayWs.Range("D89").Value = 540691

$weekWs = $wb.Worksheets.Item("week")

$weekWs.Range("A94").Value = 1
$weekWs.Range("B94").Value = "COLPAL"
$weekWs.Range("C94").Value = "Colgate Palmolive (india) Limited"
$weekWs.Range("D94").Value = "'500830"
$weekWs.Range("E94").Value = 0.88
$weekWs.Range("F94").Value = 2910.5
$weekWs.Range("G94").Value = 163914
$weekWs.Range("H94").Value = "week"
$weekWs.Range("I94").Value = "05/07/2024 11:32:27"

$weekWs.Range("A95").Value = 2
$weekWs.Range("B95").Value = "UBL"
$weekWs.Range("C95").Value = "United Breweries Limited"
$weekWs.Range("D95").Value = "'532478"
$weekWs.Range("E95").Value = 3.01
$weekWs.Range("F95").Value = 2072.6
$weekWs.Range("G95").Value = 907665
$weekWs.Range("H95").Value = "week"
$weekWs.Range("I95").Value = "05/07/2024 11:32:27"

$weekWs.Range("A96").Value = 3
$weekWs.Range("B96").Value = "NATIONALUM"
$weekWs.Range("C96").Value = "National Aluminium Company Limited"
$weekWs.Range("D96").Value = "'532234"
$weekWs.Range("E96").Value = 2.28
$weekWs.Range("F96").Value = 199.02
$weekWs.Range("G96").Value = 15659265
$weekWs.Range("H96").Value = "week"
$weekWs.Range("I96").Value = "05/07/2024 11:32:27"
